# Insert a new weekly price-report row before row 134 (Femacal de La Calera,
# Zapallo italiano), shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("134:134").Insert()

$ws.Cells.Item(134, 1).Value  = 3
$ws.Cells.Item(134, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(134, 3).Value  = "Coquimbo"
$ws.Cells.Item(134, 4).Value  = 44452
$ws.Cells.Item(134, 5).Value  = 5
$ws.Cells.Item(134, 6).Value  = 100112032
$ws.Cells.Item(134, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(134, 8).Value  = "Sin especificar"
$ws.Cells.Item(134, 9).Value  = "Primera"
$ws.Cells.Item(134, 10).Value = 190
$ws.Cells.Item(134, 11).Value = 12000
$ws.Cells.Item(134, 12).Value = 13000
$ws.Cells.Item(134, 13).Value = 12474
$ws.Cells.Item(134, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(134, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(134, 16).Value = 178
$ws.Cells.Item(134, 17).Value = 70
$ws.Cells.Item(134, 18).Value = "Hortaliza"
